# Generate Report for Handback
#
# This script mirrors the "handback" pass of the localization report:
#   - The "Status" for each tracked file flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview!E2/F2, zh-cn!C2, de-de!C2).
#   - Each language sheet's row gets its "Latest Target File" (I2) and
#     "Latest Handback File" (J2) populated, and "Latest Handback DateTime"
#     (K2) stamped with the real handback time (replacing the
#     0001-01-01 00:00:00 placeholder).
#   - A hyperlink to the source markdown file is added on I2, matching the
#     existing A2 hyperlink.
#   - The Overview/zh-cn/de-de "Status" columns and the "Latest Target
#     File"/"Latest Handback File" columns widen to fit the new text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$sourceMdFile = "d963e7ae-6d89-4111-ae40-56042b4814f9.md"
$sourceMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cad10e00046eb8b3f7a6b0b07ba478fbb0ae7a70/e2e/d963e7ae-6d89-4111-ae40-56042b4814f9.md"

# ---------------------------------------------------------------------
# Overview sheet: update the zh-cn / de-de status columns (E2 / F2)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusText

$zhCnHandoffFile = $wsZhCn.Range("G2").Value2
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $sourceMdUrl, [Type]::Missing, [Type]::Missing, $sourceMdFile) | Out-Null
$wsZhCn.Range("J2").Value = $zhCnHandoffFile
$wsZhCn.Range("K2").Value = "2016-08-16 10:59:09"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZhCn.Columns.Item(9).ColumnWidth = 39.16666666666667
$wsZhCn.Columns.Item(10).ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusText

$deDeHandoffFile = $wsDeDe.Range("G2").Value2
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $sourceMdUrl, [Type]::Missing, [Type]::Missing, $sourceMdFile) | Out-Null
$wsDeDe.Range("J2").Value = $deDeHandoffFile
$wsDeDe.Range("K2").Value = "2016-08-16 10:59:23"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDeDe.Columns.Item(9).ColumnWidth = 39.16666666666667
$wsDeDe.Columns.Item(10).ColumnWidth = 39.16666666666667
